# Insert a new price-report row at row 239 (for "Provincia de Talca"),
# pushing the existing rows 239-297 down to 240-298.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(239).Insert()

$ws.Range("A239").Value = 9
$ws.Range("B239").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C239").Value = "Metropolitana"
$ws.Range("D239").Value = 44889
$ws.Range("E239").Value = 13
$ws.Range("F239").Value = 100112026
$ws.Range("G239").Value = "Haba"
$ws.Range("H239").Value = "Sin especificar"
$ws.Range("I239").Value = "Primera"
$ws.Range("J239").Value = 80
$ws.Range("K239").Value = 15000
$ws.Range("L239").Value = 15000
$ws.Range("M239").Value = 15000
$ws.Range("N239").Value = "$/saco 25 kilos"
$ws.Range("O239").Value = "Provincia de Talca"
$ws.Range("P239").Value = 600
$ws.Range("Q239").Value = 25
$ws.Range("R239").Value = "Hortaliza"
